$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.977.72"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.827.10"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  +5.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0995"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "2.090.03"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.669"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.809.30"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "35.007.22"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "0.0₃0787"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("E27").Value = "  +4.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0551"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.703"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "92.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "1.340.00"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").Value = "2.007.06"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0669"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.30%  "
